$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds an Excel date serial that needs to move
# forward by one day (45188 -> 45189) for every data row (rows 2-74).
$range = $ws.Range("C2:C74")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
